$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 34 — the BTC "Sell" order finalizes: status flips from IN PROGRESS to
# DONE, a finalized-date timestamp is recorded, and the fee is filled in.
# ---------------------------------------------------------------------------
$ws.Range("H34").Value = "DONE"
$ws.Range("I34").Value = 42853.984583333331
$ws.Range("J34").Value = "0.00167991 USDT (0.15%) "

# ---------------------------------------------------------------------------
# Row 35 — the ETH "Buy" order is cancelled; a finalized-date note (stored as
# text, not a real date) is written into the Finalized date column.
# ---------------------------------------------------------------------------
$ws.Range("H35").Value = "CANCEL"
$ws.Range("I35").Value = " 2017-05-02 13:41:17"

# ---------------------------------------------------------------------------
# Two brand-new orders were placed: a Buy of ETC and a Sell of BTC, both
# still "IN PROGRESS". Clone the formatting of the two nearest template rows
# (33 -> 36, 34 -> 37) so number formats / wrap styles line up, then fill in
# the values.
# ---------------------------------------------------------------------------
$ws.Range("A33:I33").Copy()
$ws.Range("A36:I36").PasteSpecial(-4122)

$ws.Range("A34:I34").Copy()
$ws.Range("A37:I37").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 36 - Buy ETC
$ws.Range("A36").Value = 42857.570335648146
$ws.Range("B36").Value = "            Buy"
$ws.Range("B36").Characters(13, 3).Font.Color = 5287936
$ws.Range("C36").Value = "        ETC"
$ws.Range("D36").Value = "                 6.80`n"
$ws.Range("E36").Value = "          6.66USDT"
$ws.Range("F36").Value = "        3.4ETC"
$ws.Range("G36").Value = " ETC/USDT0000001"
$ws.Range("H36").Value = "IN PROGRESS"
$ws.Range("I36").Value = ""

# Row 37 - Sell BTC
$ws.Range("A37").Value = 42857.730462962965
$ws.Range("B37").Value = "            Sell"
$ws.Range("B37").Characters(13, 4).Font.Color = 255
$ws.Range("C37").Value = "        BTC"
$ws.Range("D37").Value = 1571
$ws.Range("E37").Value = "              1600USDT"
$ws.Range("F37").Value = "     0.00061456 BTC"
$ws.Range("G37").Value = "  BTC/USDT"
$ws.Range("H37").Value = "IN PROGRESS"
$ws.Range("I37").Value = ""

# ---------------------------------------------------------------------------
# Move the current selection to where the user ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("D41").Select()
